$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# ---------------------------------------------------------------------
# 1. Summary block (top-right mini table): component counts changed
#    because a second connector (J2) and a new resistor (R1) were added.
# ---------------------------------------------------------------------
$ws.Range("F2").Value = 4                        # Component Groups:  3 -> 4
$ws.Range("F3").Value = "27 (27 SMD/ 0 THT)"      # Component Count:   25 -> 27
$ws.Range("F4").Value = "27 (27 SMD/ 0 THT)"      # Fitted Components: 25 -> 27
$ws.Range("F6").Value = 27                        # Total Components:  25 -> 27

# ---------------------------------------------------------------------
# 2. Existing connector row (row 11) now covers J1 *and* J2, is used for
#    the "middle" position, and so its quantity doubles to 2.
# ---------------------------------------------------------------------
$ws.Range("D11").Value = "J1 J2"
$ws.Range("E11").Value = "middle"
$ws.Range("G11").Value = "'2"
# Restore G11's original cell style (quantity column keeps style index 5 -
# the same banding as A11/H11); writing the value alone would otherwise
# pick up a stray "quoted text" variant of that style.
$ws.Range("H11").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. New row 12: a 330 ohm resistor (R1) added to the BoM.
#    Build it by cloning row 10's formatting (its style alternates with
#    row 9/11) so the banding pattern keeps going, then fill in values.
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "'4"
$ws.Range("B12").Value = "Resistor"
$ws.Range("C12").Value = "R"
$ws.Range("D12").Value = "R1"
$ws.Range("E12").Value = "'330"
$ws.Range("F12").Value = "R_0201_0603Metric_Pad0.64x0.40mm_HandSolder"
$ws.Range("G12").Value = "'1"
$ws.Range("H12").Value = " "
$ws.Range("I12").Value = "~"

$ws.Range("A10:I10").Copy()
$ws.Range("A12:I12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Datasheet cell for the new resistor row is flagged (no datasheet
# available) with a light red/salmon fill.
$ws.Range("I12").Interior.Color = 9079551

# ---------------------------------------------------------------------
# 4. Column E got narrower to make room for the extra row.
# ---------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 19.7109375
